$d = $word.ActiveDocument

# 1. Update the date in the author line: "7 Mar 2023" -> "8 Mar 2023"
$d.Content.Find.Execute(" 7 Mar 2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 8 Mar 2023", 2)

# 2. Update the git revision: "a291098" -> "5eedeef"
$d.Content.Find.Execute("a291098", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5eedeef", 2)

# 3. Redraft the Methods paragraph explaining the statistical analysis.
#    Find the paragraph that begins with the old lead sentence and replace
#    its whole content (excluding the trailing paragraph mark) with the new
#    text. We use Range.Text assignment (rather than Find.Execute's
#    Replacement) so that straight apostrophes in the new text are not
#    mangled into curly ones by Word's "smart quotes" AutoCorrect, which
#    only applies to typed/Replacement text, not direct Range.Text sets.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Except as noted, all statistical analyses")) {
        $target = $p
        break
    }
}
$rng = $target.Range
$rng.End = $rng.End - 1
$rng.Text = 'All statistical analyses were performed as specified in our protocol using Stata 16 (StataCorp LLC, College Station, Texas, USA), except for one secondary analysis (see Protocol Deviations). The study is retrospective, and reviews were not randomized to use recommended ML versus no ML (for example). We therefore modelled ML use as an endogenously assigned treatment predicted by field (healthcare or welfare) and pre-specification (existence of a protocol), as planned. Resource use was analyzed using extended interval regression (Stata''s eintreg command) and time-to-completion was analyzed using a likelihood-adjusted-censoring inverse-probability-weighted regression adjustment model (LAC-IPWRA; Stata''s stteffects command). Ongoing reviews were right censored at the end of data collection (31 January 2023) and all analyses accounted for this censoring. We had no reason to suspect informative (nonrandom) censoring, so did not model a censoring mechanism. We re-expressed all estimates as ratios (relative resource use and relative time-to-completion) to aid generalization to other institutions. We did this by exponentiating differences in log resource use, and by computing ratios of mean times-to-completion using the delta method. We present two-sided 95% confidence intervals and p-values where appropriate and use a prespecified p < 0.05 significance criterion throughout. We also present the time-to-completion data using Kaplan-Meier estimates of survivor functions.'
